$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Update the "Datos actualizados..." timestamp footer (A1)
# ---------------------------------------------------------------
$ws.Cells.Item(1,1).Value2 = "Datos actualizados a 30 de Abril de 2020 a las 16:22"

# ---------------------------------------------------------------
# 2) Countries re-sorted: Moldavia overtakes Argelia (rows 58-59)
#    Swap the country names between the two rows; the per-row
#    statistics (columns B:H) are updated separately below so
#    that each country keeps / receives the right numbers.
# ---------------------------------------------------------------
$name58 = $ws.Cells.Item(58,1).Value2
$name59 = $ws.Cells.Item(59,1).Value2
$ws.Cells.Item(58,1).Value2 = $name59
$ws.Cells.Item(59,1).Value2 = $name58

# ---------------------------------------------------------------
# 3) Countries re-sorted: Maldivas jumps to the top of the
#    Jamaica..Isla de Man block (rows 118-128), pushing the rest
#    down by one position.
# ---------------------------------------------------------------
$blockRows = 118..128
$names = @()
foreach ($r in $blockRows) {
    $names += $ws.Cells.Item($r,1).Value2
}
$newNames = @($names[10]) + $names[0..9]
for ($i = 0; $i -lt $blockRows.Count; $i++) {
    $ws.Cells.Item($blockRows[$i],1).Value2 = $newNames[$i]
}

# ---------------------------------------------------------------
# 4) Refreshed statistics (Casos totales, Nuevos casos, Casos
#    activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# ---------------------------------------------------------------

# Estados Unidos
$ws.Cells.Item(4,2).Value2 = 1065956
$ws.Cells.Item(4,3).Value2 = 1762
$ws.Cells.Item(4,5).Value2 = 856687
$ws.Cells.Item(4,7).Value2 = 140
$ws.Cells.Item(4,8).Value2 = 61796

# Paises Bajos
$ws.Cells.Item(17,6).Value2 = 783

# Row 58 (now Moldavia)
$ws.Cells.Item(58,2).Value2 = 3897
$ws.Cells.Item(58,3).Value2 = 126
$ws.Cells.Item(58,4).Value2 = 1182
$ws.Cells.Item(58,5).Value2 = 2599
$ws.Cells.Item(58,6).Value2 = 237
$ws.Cells.Item(58,7).Value2 = 5
$ws.Cells.Item(58,8).Value2 = 116

# Row 59 (now Argelia)
$ws.Cells.Item(59,2).Value2 = 3848
$ws.Cells.Item(59,3).Value2 = 0
$ws.Cells.Item(59,4).Value2 = 1702
$ws.Cells.Item(59,5).Value2 = 1702
$ws.Cells.Item(59,6).Value2 = 22
$ws.Cells.Item(59,7).Value2 = 0
$ws.Cells.Item(59,8).Value2 = 444

# Uzbekistan
$ws.Cells.Item(71,4).Value2 = 1126
$ws.Cells.Item(71,5).Value2 = 882

# Row 118 (now Maldivas)
$ws.Cells.Item(118,3).Value2 = 118
$ws.Cells.Item(118,4).Value2 = 17
$ws.Cells.Item(118,5).Value2 = 378
$ws.Cells.Item(118,6).Value2 = 2
$ws.Cells.Item(118,8).Value2 = 1

# Row 119 (now Jamaica)
$ws.Cells.Item(119,3).Value2 = 0
$ws.Cells.Item(119,4).Value2 = 29
$ws.Cells.Item(119,5).Value2 = 360
$ws.Cells.Item(119,6).Value2 = 3
$ws.Cells.Item(119,7).Value2 = 0
$ws.Cells.Item(119,8).Value2 = 7

# Row 120 (now Kenia)
$ws.Cells.Item(120,2).Value2 = 396
$ws.Cells.Item(120,3).Value2 = 12
$ws.Cells.Item(120,4).Value2 = 144
$ws.Cells.Item(120,5).Value2 = 235
$ws.Cells.Item(120,6).Value2 = 2
$ws.Cells.Item(120,7).Value2 = 2
$ws.Cells.Item(120,8).Value2 = 17

# Row 121 (now El Salvador)
$ws.Cells.Item(121,2).Value2 = 395
$ws.Cells.Item(121,3).Value2 = 18
$ws.Cells.Item(121,4).Value2 = 118
$ws.Cells.Item(121,5).Value2 = 268
$ws.Cells.Item(121,6).Value2 = 3
$ws.Cells.Item(121,8).Value2 = 9

# Row 122 (now Sudan)
$ws.Cells.Item(122,2).Value2 = 375
$ws.Cells.Item(122,4).Value2 = 32
$ws.Cells.Item(122,5).Value2 = 315
$ws.Cells.Item(122,8).Value2 = 28

# Row 123 (now Estado de Palestina)
$ws.Cells.Item(123,2).Value2 = 344
$ws.Cells.Item(123,4).Value2 = 71
$ws.Cells.Item(123,5).Value2 = 271
$ws.Cells.Item(123,6).Value2 = 0
$ws.Cells.Item(123,8).Value2 = 2

# Row 124 (now Mauricio)
$ws.Cells.Item(124,2).Value2 = 332
$ws.Cells.Item(124,4).Value2 = 310
$ws.Cells.Item(124,5).Value2 = 12
$ws.Cells.Item(124,6).Value2 = 3

# Row 125 (now Venezuela)
$ws.Cells.Item(125,2).Value2 = 331
$ws.Cells.Item(125,4).Value2 = 142
$ws.Cells.Item(125,5).Value2 = 179
$ws.Cells.Item(125,8).Value2 = 10

# Row 126 (now Montenegro)
$ws.Cells.Item(126,2).Value2 = 322
$ws.Cells.Item(126,4).Value2 = 206
$ws.Cells.Item(126,5).Value2 = 109
$ws.Cells.Item(126,6).Value2 = 2
$ws.Cells.Item(126,8).Value2 = 7

# Row 127 (now Guinea Ecuatorial)
$ws.Cells.Item(127,2).Value2 = 315
$ws.Cells.Item(127,4).Value2 = 9
$ws.Cells.Item(127,5).Value2 = 305
$ws.Cells.Item(127,6).Value2 = 0
$ws.Cells.Item(127,8).Value2 = 1

# Row 128 (now Isla de Man)
$ws.Cells.Item(128,2).Value2 = 313
$ws.Cells.Item(128,3).Value2 = 0
$ws.Cells.Item(128,4).Value2 = 258
$ws.Cells.Item(128,5).Value2 = 34
$ws.Cells.Item(128,6).Value2 = 21
$ws.Cells.Item(128,8).Value2 = 21

# Suazilandia
$ws.Cells.Item(151,4).Value2 = 12
$ws.Cells.Item(151,5).Value2 = 87
